$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header / summary value updates ---
$ws.Range("E11").Value = 56940
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 1

# --- Apply the closing ("last row") border style from the old last data row (24)
#     onto row 17, which will become the sole remaining data row. ---
$ws.Range("B24:J24").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)

# --- Row 16: first worker data row gets new worker/period/value ---
$ws.Range("C16").Value = "1143337423"
$ws.Range("D16").Value = "GLADYS ESTHER MAZZEO ORTIZ"
$ws.Range("E16").Value = "2508"
$ws.Range("F16").Value = 7592

# --- Row 17: second worker data row (now the last row) ---
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "73195766"
$ws.Range("D17").Value = "TONY YAHIR CASTRO CARDOZA"
$ws.Range("E17").Value = "2508"
$ws.Range("F17").Value = 49348
$ws.Range("G17").Value = 1423500

# --- Remove the now-obsolete extra worker/period rows (18-24) ---
$ws.Rows("18:24").Delete()

# --- Column D best-fit width now that the longest name is shorter ---
$ws.Columns("D").AutoFit()

Write-Output "done"
